$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; change it to the text "1"
# (kept as text, not a number, to match the original string cell type).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"

